$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: update description text (C4) to reflect "el portal" wording ---
$ws.Range("C4").Value = "En el portal de myShopify me intento registrar sin embargo el capcha no me deja debido a que es un robot"

# --- Row 5: update description text (C5) to reflect "el portal" wording ---
$ws.Range("C5").Value = "En el portal de myshopify en el home utilizo el buscador y me aparece un resultado relacionado"

# --- Row 6: brand-new scenario "BusquedaFallida" ---
$ws.Range("B6").Value = "BusquedaFallida"
$ws.Range("C6").Value = "En el portal de myshopify utlizo el buscador con algo que no se encuentre en la tienda, debe aparecer el mensaje indicando que no hubieron resultados."
$ws.Range("D6").Value = "Me encuentro en el home "
$ws.Range("E6").Value = "1. ingresar en el buscador  una palabra de algo que no se vende en la tienda. 2. verificar que se encuentre el mensaje que le avisa al usuario que no se encontraro resultados."
$ws.Range("F6").Value = "Me aparece el mensaje 'No results found for'"

# Match formatting of row 5 (the template data row) for new row 6
$ws.Range("B5:F5").Copy()
$ws.Range("B6:F6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(6).RowHeight = 90

# --- Update active selection to C10 ---
$ws.Range("C10").Select()
